# The workbook rows 34-49 (inclusive) had their data records shuffled
# around between row positions (row 41 kept its original content).
# This script re-creates that shuffle by moving whole-row values
# (columns A:AY) between the affected rows.
#
# $targetToSource[<row>] = <row whose original content ends up there>
$targetToSource = @{
    34 = 49
    35 = 38
    36 = 46
    37 = 35
    38 = 48
    39 = 40
    40 = 43
    41 = 41
    42 = 44
    43 = 45
    44 = 34
    45 = 42
    46 = 36
    47 = 37
    48 = 39
    49 = 47
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues($row) {
    return $ws.Range("A${row}:AY${row}").Value()
}

function Set-RowValues($row, $values) {
    # Columns Y and AA hold plain-text dates (e.g. "2023-08-22"). Writing a
    # date-looking string through COM auto-converts it to a real date
    # unless the destination cell is pre-formatted as Text.
    $ws.Range("Y${row}").NumberFormat = "@"
    $ws.Range("AA${row}").NumberFormat = "@"
    $ws.Range("A${row}:AY${row}").Value = $values
    $ws.Range("Y${row}").NumberFormat = "General"
    $ws.Range("AA${row}").NumberFormat = "General"
}

# Apply the permutation using cycle decomposition so every row's original
# content is captured before it gets overwritten.
$visited = @{}

foreach ($startRow in ($targetToSource.Keys | Sort-Object)) {
    if ($visited.ContainsKey($startRow)) {
        continue
    }

    # Build the cycle starting at $startRow: startRow <- cyc1 <- cyc2 <- ... <- startRow
    $cycle = New-Object System.Collections.Generic.List[int]
    $cycle.Add($startRow) | Out-Null
    $visited[$startRow] = $true
    $cur = $targetToSource[$startRow]
    while ($cur -ne $startRow) {
        $cycle.Add($cur) | Out-Null
        $visited[$cur] = $true
        $cur = $targetToSource[$cur]
    }

    if ($cycle.Count -le 1) {
        continue  # row content unchanged
    }

    # Save the first row's original content, then shift every other row's
    # content into the previous slot, finally dropping the saved content
    # into the last slot.
    $scratch = Get-RowValues $cycle[0]
    for ($i = 0; $i -lt $cycle.Count - 1; $i++) {
        $srcValues = Get-RowValues $cycle[$i + 1]
        Set-RowValues $cycle[$i] $srcValues
    }
    Set-RowValues $cycle[$cycle.Count - 1] $scratch
}
